$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet SCD0263 -> SCD0016
$ws.Name = "SCD0016"

# Update TC_ID column (B) for rows 3-8: DGS-278 -> SCD0016-037
$ws.Range("B3").Value = "SCD0016-037"
$ws.Range("B4").Value = "SCD0016-037"
$ws.Range("B5").Value = "SCD0016-037"
$ws.Range("B6").Value = "SCD0016-037"
$ws.Range("B7").Value = "SCD0016-037"
$ws.Range("B8").Value = "SCD0016-037"

# Widen column B to fit the new, longer TC_ID text
$ws.Columns("B").ColumnWidth = 13.1

# Scroll / selection state: user ended up viewing row 8, selecting B9
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B9").Select()
